$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 442.14285
$ws.Range("I2").Value = 699
$ws.Range("J2").Value = 99.666664
$ws.Range("K2").Value = 699
$ws.Range("L2").Value = 99.666664
$ws.Range("M2").Value = -586
$ws.Range("N2").Value = -325.666664
$ws.Range("H17").Value = 11928.26
$ws.Range("J17").Value = 11928.26
$ws.Range("L17").Value = 35784.78
$ws.Range("N17").Value = -36120.78
$ws.Range("H32").Value = 5666.5
$ws.Range("I32").Value = 6124.75
$ws.Range("J32").Value = 4750
$ws.Range("K32").Value = 6124.75
$ws.Range("L32").Value = 4750
$ws.Range("M32").Value = -5798.75
$ws.Range("N32").Value = -5402
$ws.Range("H40").Value = 3999.3333
$ws.Range("J40").Value = 3999.3333
$ws.Range("L40").Value = 3999.3333
$ws.Range("N40").Value = -4349.3333
$ws.Range("H62").Value = 88857.086
$ws.Range("I62").Value = 115808.78
$ws.Range("K62").Value = 115808.78
$ws.Range("M62").Value = -115184.78
$ws.Range("H65").Value = 88857.086
$ws.Range("I65").Value = 115808.78
$ws.Range("K65").Value = 579043.9
$ws.Range("M65").Value = -575923.9
$ws.Range("H74").Value = 4758.1665
$ws.Range("I74").Value = 4025.8
$ws.Range("K74").Value = 4025.8
$ws.Range("M74").Value = -3089.8
$ws.Range("H77").Value = 4758.1665
$ws.Range("I77").Value = 4025.8
$ws.Range("K77").Value = 20129
$ws.Range("M77").Value = -15449
$ws.Range("H132").Value = 3903.7576
$ws.Range("I132").Value = 1188.45
$ws.Range("K132").Value = 3565.35
$ws.Range("M132").Value = -1035.35
$ws.Range("H135").Value = 31275.908
$ws.Range("I135").Value = 876.13635
$ws.Range("K135").Value = 7885.22715
$ws.Range("M135").Value = -5350.22715
$ws.Range("H137").Value = 49473.145
$ws.Range("I137").Value = 1966.8823
$ws.Range("J137").Value = 251374.75
$ws.Range("K137").Value = 5900.6469
$ws.Range("L137").Value = 754124.25
$ws.Range("M137").Value = -3350.6469
$ws.Range("N137").Value = -759224.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3660.5454
$ws.Range("I2").Value = 3406.125
$ws.Range("J2").Value = 4339
$ws.Range("K2").Value = 3406.125
$ws.Range("L2").Value = 4339
$ws.Range("M2").Value = -3293.125
$ws.Range("N2").Value = -4565
$ws.Range("H32").Value = 45105.855
$ws.Range("I32").Value = 25633.143
$ws.Range("K32").Value = 25633.143
$ws.Range("M32").Value = -25346.143
$ws.Range("H61").Value = 1956
$ws.Range("I61").Value = 1956
$ws.Range("K61").Value = 1956
$ws.Range("M61").Value = -1744
$ws.Range("H116").Value = 3660.5454
$ws.Range("I116").Value = 3406.125
$ws.Range("J116").Value = 4339
$ws.Range("K116").Value = 3406.125
$ws.Range("L116").Value = 4339
$ws.Range("M116").Value = -1112.125
$ws.Range("N116").Value = -8927
$ws.Range("H136").Value = 1956
$ws.Range("I136").Value = 1956
$ws.Range("K136").Value = 5868
$ws.Range("M136").Value = -3318

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3660.5454
$ws.Range("I3").Value = 3406.125
$ws.Range("J3").Value = 4339
$ws.Range("K3").Value = 3406.125
$ws.Range("L3").Value = 4339
$ws.Range("M3").Value = -3292.125
$ws.Range("N3").Value = -4567
$ws.Range("H20").Value = 7487.533
$ws.Range("I20").Value = 5124.0586
$ws.Range("K20").Value = 5124.0586
$ws.Range("M20").Value = -4877.0586
$ws.Range("H107").Value = 7527.4736
$ws.Range("I107").Value = 4273.143
$ws.Range("J107").Value = 16639.6
$ws.Range("K107").Value = 4273.143
$ws.Range("L107").Value = 16639.6
$ws.Range("M107").Value = -2353.143
$ws.Range("N107").Value = -20479.6
$ws.Range("H134").Value = 2182.2856
$ws.Range("I134").Value = 1563.75
$ws.Range("K134").Value = 4691.25
$ws.Range("M134").Value = -2156.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = 50
$ws.Range("H31").Value = 2690.7273
$ws.Range("I31").Value = 2459.8
$ws.Range("K31").Value = 2459.8
$ws.Range("M31").Value = -2164.8
$ws.Range("H34").Value = 2690.7273
$ws.Range("I34").Value = 2459.8
$ws.Range("K34").Value = 2459.8
$ws.Range("M34").Value = -2257.8
$ws.Range("H58").Value = 1926.25
$ws.Range("I58").Value = 1882.4
$ws.Range("J58").Value = 1999.3334
$ws.Range("K58").Value = 1882.4
$ws.Range("L58").Value = 1999.3334
$ws.Range("M58").Value = -1679.4
$ws.Range("N58").Value = -2405.3334
$ws.Range("H86").Value = 5864.5
$ws.Range("I86").Value = 5864.5
$ws.Range("K86").Value = 5864.5
$ws.Range("M86").Value = -4741.5
$ws.Range("H89").Value = 5864.5
$ws.Range("I89").Value = 5864.5
$ws.Range("K89").Value = 29322.5
$ws.Range("M89").Value = -23706.5
$ws.Range("H99").Value = 27879.4
$ws.Range("I99").Value = 66249.5
$ws.Range("J99").Value = 2299.3333
$ws.Range("K99").Value = 66249.5
$ws.Range("L99").Value = 2299.3333
$ws.Range("M99").Value = -64751.5
$ws.Range("N99").Value = -5295.3333
$ws.Range("H126").Value = 27879.4
$ws.Range("I126").Value = 66249.5
$ws.Range("J126").Value = 2299.3333
$ws.Range("K126").Value = 198748.5
$ws.Range("L126").Value = 6897.999899999999
$ws.Range("M126").Value = -196278.5
$ws.Range("N126").Value = -11837.9999
$ws.Range("H132").Value = 1971.92
$ws.Range("I132").Value = 1998.5416
$ws.Range("J132").Value = 1333
$ws.Range("K132").Value = 5995.6248
$ws.Range("L132").Value = 3999
$ws.Range("M132").Value = -3465.6248
$ws.Range("N132").Value = -9059
$ws.Range("H134").Value = 4450.1113
$ws.Range("I134").Value = 4864.5713
$ws.Range("K134").Value = 14593.7139
$ws.Range("M134").Value = -12058.7139
$ws.Range("H136").Value = 1926.25
$ws.Range("I136").Value = 1882.4
$ws.Range("J136").Value = 1999.3334
$ws.Range("K136").Value = 5647.200000000001
$ws.Range("L136").Value = 5998.0002
$ws.Range("M136").Value = -3097.200000000001
$ws.Range("N136").Value = -11098.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 609.36365
$ws.Range("J12").Value = 617
$ws.Range("L12").Value = 1851
$ws.Range("N12").Value = -2197
$ws.Range("H56").Value = 9896.058999999999
$ws.Range("I56").Value = 9896.058999999999
$ws.Range("K56").Value = 9896.058999999999
$ws.Range("M56").Value = -9366.058999999999
$ws.Range("H76").Value = 11000.333
$ws.Range("I76").Value = 3013
$ws.Range("J76").Value = 14994
$ws.Range("K76").Value = 9039
$ws.Range("L76").Value = 44982
$ws.Range("M76").Value = -8656
$ws.Range("N76").Value = -45748
$ws.Range("H79").Value = 11000.333
$ws.Range("I79").Value = 3013
$ws.Range("J79").Value = 14994
$ws.Range("K79").Value = 9039
$ws.Range("L79").Value = 44982
$ws.Range("M79").Value = -7713
$ws.Range("N79").Value = -47634
$ws.Range("H114").Value = 16667599
$ws.Range("I114").Value = 25000726
$ws.Range("K114").Value = 75002178
$ws.Range("M114").Value = -74998924
$ws.Range("H122").Value = 1004.9231
$ws.Range("I122").Value = 726.25
$ws.Range("J122").Value = 1076.8387
$ws.Range("K122").Value = 6536.25
$ws.Range("L122").Value = 9691.5483
$ws.Range("M122").Value = -4086.25
$ws.Range("N122").Value = -14591.5483
$ws.Range("H132").Value = 1007.35297
$ws.Range("I132").Value = 884.0909
$ws.Range("J132").Value = 1233.3334
$ws.Range("K132").Value = 7956.8181
$ws.Range("L132").Value = 11100.0006
$ws.Range("M132").Value = -5426.8181
$ws.Range("N132").Value = -16160.0006
$ws.Range("H133").Value = 10210.571
$ws.Range("I133").Value = 9094.799999999999
$ws.Range("K133").Value = 27284.4
$ws.Range("M133").Value = -22224.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6426.75
$ws.Range("I70").Value = 6235.8335
$ws.Range("J70").Value = 6999.5
$ws.Range("K70").Value = 6235.8335
$ws.Range("L70").Value = 6999.5
$ws.Range("M70").Value = -5965.8335
$ws.Range("N70").Value = -7539.5
$ws.Range("H73").Value = 6426.75
$ws.Range("I73").Value = 6235.8335
$ws.Range("J73").Value = 6999.5
$ws.Range("K73").Value = 6235.8335
$ws.Range("L73").Value = 6999.5
$ws.Range("M73").Value = -5299.8335
$ws.Range("N73").Value = -8871.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 19373.5
$ws.Range("I45").Value = 13333
$ws.Range("K45").Value = 13333
$ws.Range("M45").Value = -12926
$ws.Range("H68").Value = 1676.9333
$ws.Range("I68").Value = 1360.75
$ws.Range("K68").Value = 1360.75
$ws.Range("M68").Value = -611.75
$ws.Range("H71").Value = 1676.9333
$ws.Range("I71").Value = 1360.75
$ws.Range("K71").Value = 6803.75
$ws.Range("M71").Value = -3059.75
$ws.Range("H136").Value = 3110.3333
$ws.Range("I136").Value = 2603.7856
$ws.Range("J136").Value = 3819.5
$ws.Range("K136").Value = 7811.3568
$ws.Range("L136").Value = 11458.5
$ws.Range("M136").Value = -5261.3568
$ws.Range("N136").Value = -16558.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 12395.353
$ws.Range("I132").Value = 18254.312
$ws.Range("K132").Value = 54762.936
$ws.Range("M132").Value = -52232.936
$ws.Range("H136").Value = 864
$ws.Range("I136").Value = 891.5714
$ws.Range("J136").Value = 285
$ws.Range("K136").Value = 2674.7142
$ws.Range("L136").Value = 855
$ws.Range("M136").Value = -124.7142000000003
$ws.Range("N136").Value = -5955

Write-Output "Applied 248 cell updates across 8 sheets"